$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-dimension:* -> iaest-measure:* (curated dimensions are now measures)
$ws.Range("A2").Value = "iaest-measure:grado"
$ws.Range("D2").Value = "iaest-measure:grandes-grupos"
$ws.Range("H2").Value = "iaest-measure:sexo"

# Row 3: dim -> medida for the columns whose row-2 label became a measure
$ws.Range("A3").Value = "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Row 4: skos:Concept -> xsd:int for the same columns
$ws.Range("A4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"

# Row 5 (the mapping-*.xlsx row) is no longer produced; remove it entirely
$ws.Rows(5).Delete()
